$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Step05A")

# Row 12 - TrainedNeuralNetworkScenario10.mat scenario (small training set, ReLU)
$ws.Range("A12").Value = "TrainedNeuralNetworkScenario10.mat"
$ws.Range("B12").Value = "TrainingAndTestDataScenario3.mat"
$ws.Range("C12").Value = "NeuralNetworkScenario4.mat"
$ws.Range("F12").Value = "^"
$ws.Range("G12").Value = "^"
$ws.Range("H12").Value = 0.00015
$ws.Range("I12").Value = "^"
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = "^"
$ws.Range("L12").Value = 0.499
$ws.Range("M12").Value = 0.425
$ws.Range("N12").Value = 0.2687
$ws.Range("O12").Value = 0.3593
$ws.Range("P12").Value = "7/30/2023"
$ws.Range("Q12").Value = "small training set.  Using ReLU"

# Row 13 - TrainedNeuralNetworkScenario11.mat scenario (small training set, ReLU)
$ws.Range("A13").Value = "TrainedNeuralNetworkScenario11.mat"
$ws.Range("B13").Value = "^"
$ws.Range("C13").Value = "NeuralNetworkScenario5.mat"
$ws.Range("F13").Value = "^"
$ws.Range("G13").Value = "^"
$ws.Range("H13").Value = 0.0075
$ws.Range("I13").Value = "^"
$ws.Range("J13").Value = 600
$ws.Range("K13").Value = "^"
$ws.Range("L13").Value = 0.7885
$ws.Range("M13").Value = 0.685
$ws.Range("N13").Value = 0.1098
$ws.Range("O13").Value = 0.2112
$ws.Range("P13").Value = "7/31/2023"
$ws.Range("Q13").Value = "small training set.  Using ReLU"

# Move the active selection to A13 (matches the saved workbook state)
$ws.Range("A13").Select()
